# Timeline rework: split the single 10-row date/objective list into two
# side-by-side blocks (H:I and J:K) that fit on one printable page, reword a
# couple of objectives, and store the dates as short "mm/dd"-style text
# instead of full date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the whole working area.
$ws.Range("H1:L11").Clear()

$vCenter = -4108   # xlCenter, used for VerticalAlignment

# --- Seed new shared strings in the same order the final sheet reads them --
# (keeps the rebuilt sharedStrings table closely aligned with a natural
# left-to-right, top-to-bottom authoring pass.)
$ws.Range("I6").VerticalAlignment = $vCenter
$ws.Range("I6").Value = " tidy data, basic frames of the website"

$ws.Range("K1").NumberFormat = "@"
$ws.Range("K1").VerticalAlignment = $vCenter
$ws.Range("K1").Value = "Objective"

$leftDates = @("11/06","11/07","11/10","11/12","11/20")
for ($i = 0; $i -lt $leftDates.Length; $i++) {
  $cell = $ws.Cells.Item($i + 2, 8)   # column H
  $cell.NumberFormat = "@"
  $cell.VerticalAlignment = $vCenter
  $cell.Value = $leftDates[$i]
}

$rightDates = @("11/25","11/29","12/04","12/08","12/10")
for ($i = 0; $i -lt $rightDates.Length; $i++) {
  $cell = $ws.Cells.Item($i + 2, 10)  # column J
  $cell.NumberFormat = "@"
  $cell.VerticalAlignment = $vCenter
  $cell.Value = $rightDates[$i]
}

$ws.Range("K3").VerticalAlignment = $vCenter
$ws.Range("K3").Value = "meeting, complete report and website"

# --- Header row (re-uses "Date" / "Objectives") -----------------------------
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").VerticalAlignment = $vCenter
$ws.Range("H1").Value = "Date"
$ws.Range("I1").VerticalAlignment = $vCenter
$ws.Range("I1").Value = "Objectives"
$ws.Range("J1").NumberFormat = "@"
$ws.Range("J1").VerticalAlignment = $vCenter
$ws.Range("J1").Value = "Date"
$ws.Range("L1").NumberFormat = "@"

# --- Remaining left-block objectives (re-used text) -------------------------
$leftObjectives = @(
  "set goals and have a proposal",
  "hand in the proposal",
  "tidy the data, preliminary work",
  "meet with TA, make correction"
)
for ($i = 0; $i -lt $leftObjectives.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 9).VerticalAlignment = $vCenter
  $ws.Cells.Item($row, 9).Value = $leftObjectives[$i]
}

# --- Remaining right-block objectives (re-used text) -------------------------
$rightObjectives = @{
  2 = "started to write the final report";
  4 = "polise final report, screencast, webpages";
  5 = "do peer assessment";
  6 = "in-class discussion of projects"
}
foreach ($row in $rightObjectives.Keys) {
  $ws.Cells.Item($row, 11).VerticalAlignment = $vCenter
  $ws.Cells.Item($row, 11).Value = $rightObjectives[$row]
}

# --- View / selection / page orientation ------------------------------------
$ws.Range("H1").Select()
$ws.PageSetup.Orientation = 1
